$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of 1 Million Megaways BC slot game. Play for free and experience cascades of wins and sticky wilds in a prehistoric adventure.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Near the end of the document: drop the duplicated "Play ..." heading
#    paragraph and turn the remaining italic paragraph into the new
#    image-generation prompt text.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
$dupTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$descPara = $d.Paragraphs($count)
$descXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Create an attention-grabbing feature image for &quot;1 Million Megaways BC&quot; that captures the prehistoric theme of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be depicted holding a club and standing in front of a snowy mountain background with a caveman village in the distance. Incorporate the Mammoth and Saber-toothed tiger symbols from the game into the image, as well as the Megaways logo. Use bright colors that pop to draw in potential players and create an adventurous, exciting vibe.</w:t></w:r></w:p>'
$descPara.Range.InsertXML($descXml)
